$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Riot Platforms, Inc. (unchanged values, only MACRO_SCORE updates below)

# Row 3 - becomes Coinbase Global, Inc. / COIN
$ws.Range("B3").Value = "Coinbase Global, Inc."
$ws.Range("C3").Value = "COIN"
$ws.Range("D3").Value = 272.82
$ws.Range("E3").Value = 35.9
$ws.Range("F3").Value = 14.55
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 53
$ws.Range("J3").Value = 46
$ws.Range("K3").Value = 53

# Row 4 - becomes MARA Holdings, Inc. / MARA
$ws.Range("B4").Value = "MARA Holdings, Inc."
$ws.Range("C4").Value = "MARA"
$ws.Range("D4").Value = 11.81
$ws.Range("E4").Value = 26.1
$ws.Range("F4").Value = 15.33
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 46
$ws.Range("J4").Value = 56
$ws.Range("K4").Value = 50.2

# Row 5 - becomes Bitcoin USD / BTC-USD
$ws.Range("B5").Value = "Bitcoin USD"
$ws.Range("C5").Value = "BTC-USD"
$ws.Range("D5").Value = 85851.67999999999
$ws.Range("E5").Value = 33.5
$ws.Range("F5").Value = -1.71
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = 36
$ws.Range("I5").Value = 43
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 49

# Row 6 - Strategy Inc (unchanged values, only MACRO_SCORE updates below)

# Update MACRO_SCORE (column N) for all data rows
$ws.Range("N2").Value = 85.92500513438651
$ws.Range("N3").Value = 85.92500513438651
$ws.Range("N4").Value = 85.92500513438651
$ws.Range("N5").Value = 85.92500513438651
$ws.Range("N6").Value = 85.92500513438651
